$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.131.98"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "2.551.90"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "585.53"
$ws.Range("E5").Value = "  +2.18%  "
$ws.Range("D6").Value = "147.58"
$ws.Range("E6").Value = "  -2.32%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").Value = "5.55"
$ws.Range("E10").Value = "  -3.72%  "
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D14").Value = "3.003.60"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").Value = "63.018.11"
$ws.Range("E15").Value = "  -0.58%  "
$ws.Range("D16").Value = "0.0000144"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Value = "2.562.80"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("E18").Value = "  -2.80%  "
$ws.Range("D19").Value = "337.31"
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "65.93"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("B26").Value = "SuiNetwork"
$ws.Range("C26").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D26").Value = "1.49"
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("D28").Value = "8.42"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("D29").Value = "7.69"
$ws.Range("E29").Value = "  +8.18%  "
$ws.Range("E30").Value = "  +5.46%  "
$ws.Range("E31").Value = "  -2.19%  "
$ws.Range("D32").Value = "178.31"
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("D33").Value = "1.56"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("D34").Value = "419.68"
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("D35").Value = "19.19"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "4.38"
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "39.68"
$ws.Range("E41").Value = "  -0.78%  "
$ws.Range("D42").Value = "150.57"
$ws.Range("E42").Value = "  -2.81%  "
$ws.Range("D43").Value = "3.81"
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").Value = "20.89"
$ws.Range("E44").Value = "  -1.61%  "
$ws.Range("E45").Value = "  +2.03%  "
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("D48").Value = "0.0238"
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("E49").Value = "  -1.91%  "
$ws.Range("D50").Value = "1.73"
$ws.Range("E50").Value = "  -5.66%  "
$ws.Range("E51").Value = "  -0.37%  "
